$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 121: flight #120
$ws.Cells.Item(121, 1).Value = 120
$ws.Cells.Item(121, 2).Value = 'Friday, Jan 13'
$ws.Cells.Item(121, 3).Value = '6:10 AM'
$ws.Cells.Item(121, 4).Value = 'FR1970'
$ws.Cells.Item(121, 5).Value = 'Catania'
$ws.Cells.Item(121, 6).Value = '(CTA)'
$ws.Cells.Item(121, 7).Value = 'Ryanair '
$ws.Cells.Item(121, 8).Value = 'B738'
$ws.Cells.Item(121, 9).Value = '(SP-RKD)'
$ws.Cells.Item(121, 10).Value = '6:11 AM'
$ws.Cells.Item(121, 12).Value = '0 hours, 1 minutes'

# Row 122: flight #121
$ws.Cells.Item(122, 1).Value = 121
$ws.Cells.Item(122, 2).Value = 'Friday, Jan 13'
$ws.Cells.Item(122, 3).Value = '6:10 AM'
$ws.Cells.Item(122, 4).Value = 'FR4528'
$ws.Cells.Item(122, 5).Value = 'Oslo'
$ws.Cells.Item(122, 6).Value = '(TRF)'
$ws.Cells.Item(122, 7).Value = 'Ryanair '
$ws.Cells.Item(122, 8).Value = 'B738'
$ws.Cells.Item(122, 9).Value = '(SP-RKW)'
$ws.Cells.Item(122, 10).Value = '6:13 AM'
$ws.Cells.Item(122, 12).Value = '0 hours, 3 minutes'

# Row 123: flight #122
$ws.Cells.Item(123, 1).Value = 122
$ws.Cells.Item(123, 2).Value = 'Friday, Jan 13'
$ws.Cells.Item(123, 3).Value = '6:15 AM'
$ws.Cells.Item(123, 4).Value = 'FR2022'
$ws.Cells.Item(123, 5).Value = 'Dublin'
$ws.Cells.Item(123, 6).Value = '(DUB)'
$ws.Cells.Item(123, 7).Value = 'Buzz '
$ws.Cells.Item(123, 8).Value = 'B38M'
$ws.Cells.Item(123, 9).Value = '(SP-RZE)'
$ws.Cells.Item(123, 10).Value = '6:21 AM'
$ws.Cells.Item(123, 12).Value = '0 hours, 6 minutes'

# Row 124: flight #123
$ws.Cells.Item(124, 1).Value = 123
$ws.Cells.Item(124, 2).Value = 'Friday, Jan 13'
$ws.Cells.Item(124, 3).Value = '6:45 AM'
$ws.Cells.Item(124, 4).Value = 'FR2263'
$ws.Cells.Item(124, 5).Value = 'Lisbon'
$ws.Cells.Item(124, 6).Value = '(LIS)'
$ws.Cells.Item(124, 7).Value = 'Ryanair '
$ws.Cells.Item(124, 8).Value = 'B738'
$ws.Cells.Item(124, 9).Value = '(SP-RKP)'
$ws.Cells.Item(124, 10).Value = '7:20 AM'
$ws.Cells.Item(124, 12).Value = '0 hours, 35 minutes'

# Row 125: flight #124
$ws.Cells.Item(125, 1).Value = 124
$ws.Cells.Item(125, 2).Value = 'Friday, Jan 13'
$ws.Cells.Item(125, 3).Value = '7:15 AM'
$ws.Cells.Item(125, 4).Value = 'FR1888'
$ws.Cells.Item(125, 5).Value = 'Paris'
$ws.Cells.Item(125, 6).Value = '(BVA)'
$ws.Cells.Item(125, 7).Value = 'Buzz '
$ws.Cells.Item(125, 8).Value = 'B38M'
$ws.Cells.Item(125, 9).Value = '(SP-RZG)'
$ws.Cells.Item(125, 10).Value = '7:11 AM'
$ws.Cells.Item(125, 12).Value = '0 hours, -4 minutes'

# Row 126: flight #125
$ws.Cells.Item(126, 1).Value = 125
$ws.Cells.Item(126, 2).Value = 'Friday, Jan 13'
$ws.Cells.Item(126, 3).Value = '7:20 AM'
$ws.Cells.Item(126, 4).Value = 'FR1056'
$ws.Cells.Item(126, 5).Value = 'Brussels'
$ws.Cells.Item(126, 6).Value = '(CRL)'
$ws.Cells.Item(126, 7).Value = 'Ryanair '
$ws.Cells.Item(126, 8).Value = 'B738'
$ws.Cells.Item(126, 9).Value = '(SP-RSP)'
$ws.Cells.Item(126, 10).Value = '7:27 AM'
$ws.Cells.Item(126, 12).Value = '0 hours, 7 minutes'

# Row 127: flight #126
$ws.Cells.Item(127, 1).Value = 126
$ws.Cells.Item(127, 2).Value = 'Friday, Jan 13'
$ws.Cells.Item(127, 3).Value = '7:30 AM'
$ws.Cells.Item(127, 4).Value = 'FR3284'
$ws.Cells.Item(127, 5).Value = 'Riga'
$ws.Cells.Item(127, 6).Value = '(RIX)'
$ws.Cells.Item(127, 7).Value = 'Ryanair '
$ws.Cells.Item(127, 8).Value = 'B738'
$ws.Cells.Item(127, 9).Value = '(SP-RKT)'
$ws.Cells.Item(127, 10).Value = '7:35 AM'
$ws.Cells.Item(127, 12).Value = '0 hours, 5 minutes'

# Row 128: flight #127
$ws.Cells.Item(128, 1).Value = 127
$ws.Cells.Item(128, 2).Value = 'Friday, Jan 13'
$ws.Cells.Item(128, 3).Value = '8:20 AM'
$ws.Cells.Item(128, 4).Value = 'FR1574'
$ws.Cells.Item(128, 5).Value = 'Vienna'
$ws.Cells.Item(128, 6).Value = '(VIE)'
$ws.Cells.Item(128, 7).Value = 'Ryanair '
$ws.Cells.Item(128, 8).Value = 'B738'
$ws.Cells.Item(128, 9).Value = '(SP-RKV)'
$ws.Cells.Item(128, 10).Value = '8:23 AM'
$ws.Cells.Item(128, 12).Value = '0 hours, 3 minutes'

# Row 129: flight #128
$ws.Cells.Item(129, 1).Value = 128
$ws.Cells.Item(129, 2).Value = 'Friday, Jan 13'
$ws.Cells.Item(129, 3).Value = '9:30 AM'
$ws.Cells.Item(129, 4).Value = 'FR2008'
$ws.Cells.Item(129, 5).Value = 'London'
$ws.Cells.Item(129, 6).Value = '(STN)'
$ws.Cells.Item(129, 7).Value = 'Ryanair '
$ws.Cells.Item(129, 8).Value = 'B738'
$ws.Cells.Item(129, 9).Value = '(EI-EVP)'
$ws.Cells.Item(129, 10).Value = '9:32 AM'
$ws.Cells.Item(129, 12).Value = '0 hours, 2 minutes'

# Row 130: flight #129
$ws.Cells.Item(130, 1).Value = 129
$ws.Cells.Item(130, 2).Value = 'Friday, Jan 13'
$ws.Cells.Item(130, 3).Value = '9:45 AM'
$ws.Cells.Item(130, 4).Value = 'FR1932'
$ws.Cells.Item(130, 5).Value = 'Leeds'
$ws.Cells.Item(130, 6).Value = '(LBA)'
$ws.Cells.Item(130, 7).Value = 'Ryanair '
$ws.Cells.Item(130, 8).Value = 'B738'
$ws.Cells.Item(130, 9).Value = '(EI-EKY)'
$ws.Cells.Item(130, 10).Value = '10:33 AM'
$ws.Cells.Item(130, 12).Value = '0 hours, 48 minutes'

# Row 131: flight #130
$ws.Cells.Item(131, 1).Value = 130
$ws.Cells.Item(131, 2).Value = 'Friday, Jan 13'
$ws.Cells.Item(131, 3).Value = '10:05 AM'
$ws.Cells.Item(131, 4).Value = 'FR6945'
$ws.Cells.Item(131, 5).Value = 'Barcelona'
$ws.Cells.Item(131, 6).Value = '(BCN)'
$ws.Cells.Item(131, 7).Value = 'Ryanair '
$ws.Cells.Item(131, 8).Value = 'B738'
$ws.Cells.Item(131, 9).Value = '(EI-DYC)'
$ws.Cells.Item(131, 10).Value = '10:11 AM'
$ws.Cells.Item(131, 12).Value = '0 hours, 6 minutes'

# Row 132: flight #131
$ws.Cells.Item(132, 1).Value = 131
$ws.Cells.Item(132, 2).Value = 'Friday, Jan 13'
$ws.Cells.Item(132, 3).Value = '10:45 AM'
$ws.Cells.Item(132, 4).Value = 'FR1934'
$ws.Cells.Item(132, 5).Value = 'Helsinki'
$ws.Cells.Item(132, 6).Value = '(HEL)'
$ws.Cells.Item(132, 7).Value = 'Ryanair '
$ws.Cells.Item(132, 8).Value = 'B738'
$ws.Cells.Item(132, 9).Value = '(SP-RKW)'
$ws.Cells.Item(132, 10).Value = '10:46 AM'
$ws.Cells.Item(132, 12).Value = '0 hours, 1 minutes'

# Row 133: flight #132
$ws.Cells.Item(133, 1).Value = 132
$ws.Cells.Item(133, 2).Value = 'Friday, Jan 13'
$ws.Cells.Item(133, 3).Value = '12:25 PM'
$ws.Cells.Item(133, 4).Value = 'FR6121'
$ws.Cells.Item(133, 5).Value = 'Tenerife'
$ws.Cells.Item(133, 6).Value = '(TFS)'
$ws.Cells.Item(133, 7).Value = 'Ryanair '
$ws.Cells.Item(133, 8).Value = 'B738'
$ws.Cells.Item(133, 9).Value = '(SP-RSP)'
$ws.Cells.Item(133, 10).Value = '12:36 PM'
$ws.Cells.Item(133, 12).Value = '0 hours, 11 minutes'

# Row 134: flight #133
$ws.Cells.Item(134, 1).Value = 133
$ws.Cells.Item(134, 2).Value = 'Friday, Jan 13'
$ws.Cells.Item(134, 3).Value = '12:50 PM'
$ws.Cells.Item(134, 4).Value = 'FR4533'
$ws.Cells.Item(134, 5).Value = 'Porto'
$ws.Cells.Item(134, 6).Value = '(OPO)'
$ws.Cells.Item(134, 7).Value = 'Buzz '
$ws.Cells.Item(134, 8).Value = 'B38M'
$ws.Cells.Item(134, 9).Value = '(SP-RZG)'
$ws.Cells.Item(134, 10).Value = '12:58 PM'
$ws.Cells.Item(134, 12).Value = '0 hours, 8 minutes'

# Row 135: flight #134
$ws.Cells.Item(135, 1).Value = 134
$ws.Cells.Item(135, 2).Value = 'Friday, Jan 13'
$ws.Cells.Item(135, 3).Value = '1:20 PM'
$ws.Cells.Item(135, 4).Value = 'FR1942'
$ws.Cells.Item(135, 5).Value = 'Bologna'
$ws.Cells.Item(135, 6).Value = '(BLQ)'
$ws.Cells.Item(135, 7).Value = 'Ryanair '
$ws.Cells.Item(135, 8).Value = 'B738'
$ws.Cells.Item(135, 9).Value = '(SP-RKD)'
$ws.Cells.Item(135, 10).Value = '1:23 PM'
$ws.Cells.Item(135, 12).Value = '0 hours, 3 minutes'

# Row 136: flight #135
$ws.Cells.Item(136, 1).Value = 135
$ws.Cells.Item(136, 2).Value = 'Friday, Jan 13'
$ws.Cells.Item(136, 3).Value = '1:30 PM'
$ws.Cells.Item(136, 4).Value = 'FR1106'
$ws.Cells.Item(136, 5).Value = 'Alicante'
$ws.Cells.Item(136, 6).Value = '(ALC)'
$ws.Cells.Item(136, 7).Value = 'Ryanair '
$ws.Cells.Item(136, 8).Value = 'B738'
$ws.Cells.Item(136, 9).Value = '(SP-RKV)'
$ws.Cells.Item(136, 10).Value = '1:43 PM'
$ws.Cells.Item(136, 12).Value = '0 hours, 13 minutes'

# Row 137: flight #136
$ws.Cells.Item(137, 1).Value = 136
$ws.Cells.Item(137, 2).Value = 'Friday, Jan 13'
$ws.Cells.Item(137, 3).Value = '1:35 PM'
$ws.Cells.Item(137, 4).Value = 'FR1944'
$ws.Cells.Item(137, 5).Value = 'Stockholm'
$ws.Cells.Item(137, 6).Value = '(ARN)'
$ws.Cells.Item(137, 7).Value = 'Buzz '
$ws.Cells.Item(137, 8).Value = 'B38M'
$ws.Cells.Item(137, 9).Value = '(SP-RZE)'
$ws.Cells.Item(137, 10).Value = '1:35 PM'
$ws.Cells.Item(137, 12).Value = '0 hours, 0 minutes'
